$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '24.596.85'
$ws.Range("E2").Value = '  +3.51%  '
$ws.Range("D3").Value = '1.695.95'
$ws.Range("E3").Value = '  +1.93%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.22%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '317.69'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.51%  '
$ws.Range("E6").Value = '  +0.18%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3954'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.47%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4017'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.39%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.536'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +7.89%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.003'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.26%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '53.69'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +8.23%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08779'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.33%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.262'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +8.68%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '23.24'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.15%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.00001320'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.23%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.615'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.93%  '
$ws.Range("D17").Value = '1.697.50'
$ws.Range("E17").Value = '  +1.98%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '101.05'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.92%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.07014'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.50%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '19.70'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.86%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.892'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.30%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.001'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.18%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '14.09'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.56%  '
$ws.Range("D24").Value = '24.591.05'
$ws.Range("E24").Value = '  +3.53%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.052'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +8.45%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.340'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.85%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.36'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.46%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '159.68'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.14%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.237'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.35%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '134.47'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.60%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.474'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +14.98%  '
$ws.Range("D32").Value = '1.882.33'
$ws.Range("E32").Value = '  +1.90%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.099'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.20%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '7.444'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +13.34%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.08518'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.07%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '11.38'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +9.99%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.968'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.53%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2736'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.79%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '14.57'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.00%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.02772'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +9.25%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.09018'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.60%  '
$ws.Range("E42").Value = '  +0.83%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.7710'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.94%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.7223'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.52%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '15.42'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.37%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.528'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +5.31%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.229'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.28%  '
$ws.Range("B48").Value = 'Frax'
$ws.Range("C48").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.001'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.17%  '
$ws.Range("B49").Value = 'Flow'
$ws.Range("C49").Value = 'https://coinranking.com/coin/QQ0NCmjVq+flow-flow'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.350'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +12.75%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '141.04'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.65%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.08043'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.46%  '
